# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to market-price / profit columns (H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 7694189.5
$ws.Range("J17").Value = 7694189.5
$ws.Range("L17").Value = 23082568.5
$ws.Range("N17").Value = -23082904.5

$ws.Range("H19").Value = 1108
$ws.Range("I19").Value = 1255.8889
$ws.Range("J19").Value = 886.1667
$ws.Range("K19").Value = 1255.8889
$ws.Range("L19").Value = 886.1667
$ws.Range("M19").Value = -1080.8889
$ws.Range("N19").Value = -1236.1667

$ws.Range("H34").Value = 5349.75
$ws.Range("I34").Value = 5349.75
$ws.Range("K34").Value = 5349.75
$ws.Range("M34").Value = -5146.75

$ws.Range("H36").Value = 5349.75
$ws.Range("I36").Value = 5349.75
$ws.Range("K36").Value = 5349.75
$ws.Range("M36").Value = -4634.75

$ws.Range("H40").Value = 4900.0713
$ws.Range("I40").Value = 3108.5454
$ws.Range("J40").Value = 6059.294
$ws.Range("K40").Value = 3108.5454
$ws.Range("L40").Value = 6059.294
$ws.Range("M40").Value = -2933.5454
$ws.Range("N40").Value = -6409.294

$ws.Range("H106").Value = 4944.25
$ws.Range("I106").Value = 4771.2607
$ws.Range("K106").Value = 4771.2607
$ws.Range("M106").Value = -4140.2607

$ws.Range("H112").Value = 1912.2
$ws.Range("J112").Value = 2020.5
$ws.Range("L112").Value = 6061.5
$ws.Range("N112").Value = -8277.5

$ws.Range("H132").Value = 1355.4572
$ws.Range("I132").Value = 1052.9286
$ws.Range("K132").Value = 3158.7858
$ws.Range("M132").Value = -628.7857999999997

$ws.Range("H135").Value = 2466.1
$ws.Range("I135").Value = 2110.5
$ws.Range("J135").Value = 2999.5
$ws.Range("K135").Value = 18994.5
$ws.Range("L135").Value = 26995.5
$ws.Range("M135").Value = -16459.5
$ws.Range("N135").Value = -32065.5

$ws.Range("H137").Value = 37704.895
$ws.Range("I137").Value = 43087.062
$ws.Range("K137").Value = 129261.186
$ws.Range("M137").Value = -126711.186

$ws.Range("H138").Value = 6743.769
$ws.Range("J138").Value = 6587.2
$ws.Range("L138").Value = 19761.6
$ws.Range("N138").Value = -30041.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2926.3857
$ws.Range("I32").Value = 2182.3936
$ws.Range("J32").Value = 7969
$ws.Range("K32").Value = 2182.3936
$ws.Range("L32").Value = 7969
$ws.Range("M32").Value = -1895.3936
$ws.Range("N32").Value = -8543

$ws.Range("H61").Value = 4635446
$ws.Range("I61").Value = 6414914
$ws.Range("J61").Value = 8829.200000000001
$ws.Range("K61").Value = 6414914
$ws.Range("L61").Value = 8829.200000000001
$ws.Range("M61").Value = -6414702
$ws.Range("N61").Value = -9253.200000000001

$ws.Range("H74").Value = 26922
$ws.Range("I74").Value = 2747.2173
$ws.Range("J74").Value = 138126
$ws.Range("K74").Value = 2747.2173
$ws.Range("L74").Value = 138126
$ws.Range("M74").Value = -1873.2173
$ws.Range("N74").Value = -139874

$ws.Range("H77").Value = 26922
$ws.Range("I77").Value = 2747.2173
$ws.Range("J77").Value = 138126
$ws.Range("K77").Value = 13736.0865
$ws.Range("L77").Value = 690630
$ws.Range("M77").Value = -9368.086499999999
$ws.Range("N77").Value = -699366

$ws.Range("H110").Value = 6613.9414
$ws.Range("I110").Value = 4123.9165
$ws.Range("K110").Value = 4123.9165
$ws.Range("M110").Value = -2078.9165

$ws.Range("H132").Value = 4022.262
$ws.Range("I132").Value = 3118.879
$ws.Range("J132").Value = 7334.6665
$ws.Range("K132").Value = 9356.636999999999
$ws.Range("L132").Value = 22003.9995
$ws.Range("M132").Value = -6826.636999999999
$ws.Range("N132").Value = -27063.9995

$ws.Range("H135").Value = 130749.75
$ws.Range("J135").Value = 130749.75
$ws.Range("L135").Value = 130749.75
$ws.Range("N135").Value = -140889.75

$ws.Range("H136").Value = 4635446
$ws.Range("I136").Value = 6414914
$ws.Range("J136").Value = 8829.200000000001
$ws.Range("K136").Value = 19244742
$ws.Range("L136").Value = 26487.6
$ws.Range("M136").Value = -19242192
$ws.Range("N136").Value = -31587.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 3315.6155
$ws.Range("I22").Value = 3294
$ws.Range("K22").Value = 3294
$ws.Range("M22").Value = -3121

$ws.Range("H134").Value = 4325.5312
$ws.Range("I134").Value = 4289.933
$ws.Range("K134").Value = 12869.799
$ws.Range("M134").Value = -10334.799

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 9581.799999999999

$ws.Range("H51").Value = 19200
$ws.Range("J51").Value = 19200
$ws.Range("L51").Value = 19200
$ws.Range("N51").Value = -20672

$ws.Range("H59").Value = 29750
$ws.Range("J59").Value = 44500
$ws.Range("L59").Value = 44500
$ws.Range("N59").Value = -46790

$ws.Range("H60").Value = 14700
$ws.Range("J60").Value = 26000
$ws.Range("L60").Value = 26000
$ws.Range("N60").Value = -27022

$ws.Range("H61").Value = 19200
$ws.Range("J61").Value = 19200
$ws.Range("L61").Value = 19200
$ws.Range("N61").Value = -19896

$ws.Range("H107").Value = 2545.7144
$ws.Range("I107").Value = 3825.8572
$ws.Range("J107").Value = 1905.6428
$ws.Range("K107").Value = 3825.8572
$ws.Range("L107").Value = 1905.6428
$ws.Range("M107").Value = -1905.8572
$ws.Range("N107").Value = -5745.6428

$ws.Range("H132").Value = 6174.1577
$ws.Range("I132").Value = 5750.6
$ws.Range("K132").Value = 17251.8
$ws.Range("M132").Value = -14721.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 197.375
$ws.Range("J23").Value = 154.83333
$ws.Range("L23").Value = 464.49999
$ws.Range("N23").Value = -934.49999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 37975.9
$ws.Range("I93").Value = 36000
$ws.Range("J93").Value = 38195.445
$ws.Range("K93").Value = 36000
$ws.Range("L93").Value = 38195.445
$ws.Range("M93").Value = -34128
$ws.Range("N93").Value = -41939.445

$ws.Range("H126").Value = 8447.77
$ws.Range("I126").Value = 7315.375
$ws.Range("J126").Value = 10259.6
$ws.Range("K126").Value = 21946.125
$ws.Range("L126").Value = 30778.8
$ws.Range("M126").Value = -19476.125
$ws.Range("N126").Value = -35718.8

$ws.Range("H132").Value = 1836.5834
$ws.Range("I132").Value = 1267.375
$ws.Range("K132").Value = 3802.125
$ws.Range("M132").Value = -1272.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2282.6667
$ws.Range("I22").Value = 850
$ws.Range("K22").Value = 850
$ws.Range("M22").Value = -555

$ws.Range("H27").Value = 2282.6667
$ws.Range("I27").Value = 850
$ws.Range("K27").Value = 850
$ws.Range("M27").Value = -743

$ws.Range("H100").Value = 3574730.5
$ws.Range("I100").Value = 7145753
$ws.Range("K100").Value = 7145753
$ws.Range("M100").Value = -7145212

$ws.Range("H106").Value = 10320
$ws.Range("J106").Value = 10320
$ws.Range("L106").Value = 10320
$ws.Range("N106").Value = -12844

$ws.Range("H132").Value = 8531.357
$ws.Range("I132").Value = 9536.666999999999
$ws.Range("K132").Value = 28610.001
$ws.Range("M132").Value = -26080.001

$ws.Range("H136").Value = 2513.8076
$ws.Range("I136").Value = 2192.1738
$ws.Range("K136").Value = 6576.5214
$ws.Range("M136").Value = -4026.5214

$ws.Range("H139").Value = 80714.5
$ws.Range("J139").Value = 80714.5
$ws.Range("L139").Value = 80714.5
$ws.Range("N139").Value = -90994.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("N13").ClearContents()

$ws.Range("H100").Value = 725.5294
$ws.Range("I100").Value = 667.38464
$ws.Range("K100").Value = 1334.76928
$ws.Range("M100").Value = -793.76928

$ws.Range("H125").Value = 65000
$ws.Range("J125").Value = 65000
$ws.Range("L125").Value = 65000
$ws.Range("N125").Value = -74840

$ws.Range("H136").Value = 3332.56
$ws.Range("I136").Value = 2746.4285
$ws.Range("K136").Value = 8239.2855
$ws.Range("M136").Value = -5689.2855

$ws.Range("H139").Value = 80666.664
$ws.Range("J139").Value = 80666.664
$ws.Range("L139").Value = 80666.664
$ws.Range("N139").Value = -90946.664
